$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.7
$ws.Range("I3").Value = 4.75
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 4.33
$ws.Range("N3").Value = 1.67
$ws.Range("O3").Value = 2.15
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.2
$ws.Range("Z3").Value = 15

# Row 9 updates
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 1.53
$ws.Range("M9").Value = 3.6
$ws.Range("T9").Value = 16
$ws.Range("U9").Value = 35
$ws.Range("V9").Value = 17.5
$ws.Range("X9").Value = 55
$ws.Range("Y9").Value = 50
$ws.Range("Z9").Value = 12.5
$ws.Range("AB9").Value = 15
$ws.Range("AC9").Value = 60
$ws.Range("AD9").Value = 450
$ws.Range("AE9").Value = 8
$ws.Range("AF9").Value = 8
$ws.Range("AI9").Value = 11.5
